$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 341. This shifts the existing rows 341-430
# down to 342-431 (matching the rest of the diff, which is just every
# subsequent row's data moving down by one row).
$ws.Rows.Item(341).Insert()

# Populate the newly inserted row 341 with its data. Columns not listed in
# the diff as changed keep the same values the (old) row 341 had, which
# Excel's Insert() already preserved via the copied formatting; we set
# them explicitly here to be sure the values are correct.
$ws.Cells.Item(341, 1).Value = 4
$ws.Cells.Item(341, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(341, 3).Value = "Los Lagos"
$ws.Cells.Item(341, 4).Value = 44508
$ws.Cells.Item(341, 5).Value = 10
$ws.Cells.Item(341, 6).Value = 100112020
$ws.Cells.Item(341, 7).Value = "Tomate"
$ws.Cells.Item(341, 8).Value = "Larga vida"
$ws.Cells.Item(341, 9).Value = "Primera"
$ws.Cells.Item(341, 10).Value = 300
$ws.Cells.Item(341, 11).Value = 23000
$ws.Cells.Item(341, 12).Value = 23000
$ws.Cells.Item(341, 13).Value = 23000
$ws.Cells.Item(341, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(341, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(341, 16).Value = 1278
$ws.Cells.Item(341, 17).Value = 18
$ws.Cells.Item(341, 18).Value = "Hortaliza"
